$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the empty bold paragraphs that sit right before:
#      - "Process"          (after the title)
#      - "The Guidelines"
#      - "Tips and Tricks"
#    Each of these paragraphs has no runs/text; deleting their Range
#    (which is just the paragraph mark) merges them away, leaving the
#    heading paragraph that follows with its own, already-matching
#    pPr/rPr formatting.
#
#    We walk paragraphs from the end backwards so indices found before
#    each deletion stay valid for the ones still to come.
# ------------------------------------------------------------------
$targets = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Length -eq 1) {
        $t = $p.Range.Text
        if ($t -eq [char]13) {
            $targets += $i
        }
    }
}

# Only the 3 empty headline-spacer paragraphs (before "Process",
# "The Guidelines" and "Tips and Tricks") are removed; the empty
# paragraph further down (before "After you are done writing...")
# is left untouched.
$toDelete = @()
foreach ($idx in $targets) {
    $next = $d.Paragraphs.Item($idx + 1)
    $nextText = $next.Range.Text
    if ($nextText.StartsWith("Process") -or $nextText.StartsWith("The Guidelines") -or $nextText.StartsWith("Tips and Tricks")) {
        $toDelete += $idx
    }
}

$n = $toDelete.Count
for ($k = $n - 1; $k -ge 0; $k--) {
    $idx = $toDelete[$k]
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
}

# ------------------------------------------------------------------
# 2) Split the "The following tips and tricks ... documentation:"
#    run right after "docume" and drop a "_GoBack" bookmark at the
#    split point (this is what Word leaves behind when the cursor
#    was there at save time). "_GoBack" is a singleton bookmark, so
#    re-adding it here automatically removes it from its old spot
#    further down in the document (in front of "After you are done
#    writing...") - no separate cleanup call is required.
# ------------------------------------------------------------------
$needle = "The following tips and tricks might help you improve your documentation:"
$splitAfter = "The following tips and tricks might help you improve your docume"

$rng = $d.Content
$rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$paraStart = $rng.Start
$splitPos = $paraStart + $splitAfter.Length
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ------------------------------------------------------------------
# 3) Drop the stale <w:lastRenderedPageBreak/> that sits in front of
#    "After you are done writing, run a check ...". Re-typing the
#    run's text (via a throwaway placeholder swap) forces the host to
#    rebuild the run without the stale page-break cache marker.
# ------------------------------------------------------------------
$needle2 = "After you are done writing, run a check for the following and self-edit:"
$rng2 = $d.Content
$rng2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start2 = $rng2.Start
$textRange = $d.Range($start2, $rng2.End)
$textRange.Text = "GoBackPlaceholder"
$textRange2 = $d.Range($start2, $start2 + 17)
$textRange2.Text = $needle2

# ------------------------------------------------------------------
# 4) Shrink the page margins to 0.5" (720 twips == 36 points) on
#    every side; header/footer distances stay as-is (720 twips already).
#    PageSetup margins are expressed in points, so 720 twips / 20 = 36.
# ------------------------------------------------------------------
$d.PageSetup.TopMargin = 36
$d.PageSetup.BottomMargin = 36
$d.PageSetup.LeftMargin = 36
$d.PageSetup.RightMargin = 36
